$d = $word.ActiveDocument

# 1. Remove the old _GoBack bookmark (it will be re-created later at the new edit location)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Locate the placeholder paragraphs "ddd" and "dd" that get replaced/expanded
$idxDdd = -1
$idxDd = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "ddd") { $idxDdd = $i }
    if ($idxDdd -ne -1 -and $t -eq "dd" -and $idxDd -eq -1 -and $i -eq ($idxDdd + 1)) { $idxDd = $i }
}

if ($idxDdd -eq -1 -or $idxDd -eq -1) {
    throw "Could not locate 'ddd'/'dd' placeholder paragraphs (found ddd=$idxDdd dd=$idxDd)"
}

$pDdd = $d.Paragraphs.Item($idxDdd)
$pDd = $d.Paragraphs.Item($idxDd)
$rng = $d.Range($pDdd.Range.Start, $pDd.Range.End)
$rng.Delete()

$insertAt = $d.Paragraphs.Item($idxDdd).Range
$insertAt.Collapse(1)

$newXml = @'
<pkg:xmlData xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:p>
  <w:pPr>
    <w:jc w:val="center"/>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>sh</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:t>aders</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:t xml:space="preserve"> in external files</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Needs several includes:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="1134"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="808080"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>#include</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>&lt;</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>iostream</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>&gt;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="1134"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="808080"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>#include</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>&lt;</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>fstream</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>&gt;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:left="1134"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="808080"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>#include</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>&lt;string&gt;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Following function takes file link and </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>returnes</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> string:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>std</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>::</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="2B91AF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>string</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>SringFromFile</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="0000FF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>const</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>std</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>::</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="2B91AF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>string</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve">&amp; </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="808080"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>file_link</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>{</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>std</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>::</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="2B91AF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>ifstream</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> ifs(</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="808080"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>file_link</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>);</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>std</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>::</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="2B91AF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>string</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> s;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>std</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>::</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>getline</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>ifs, s, (</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="0000FF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>char</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>ifs.eof</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>());</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:tab/>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="0000FF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>return</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> s;</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:left="851"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>}</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>“</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>shader</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve">” variable line is changed to: </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:ind w:left="851"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="0000FF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>unsigned</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="0000FF"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>int</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>shader</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve"> = </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>CreateShader</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>SringFromFile</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>"res/</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>vertexShader.shader</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>"</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t xml:space="preserve">), </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>SringFromFile</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>"res/</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>fragmentShader.shader</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="A31515"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>"</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
    <w:t>));</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:autoSpaceDE w:val="0"/>
    <w:autoSpaceDN w:val="0"/>
    <w:adjustRightInd w:val="0"/>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
      <w:color w:val="000000"/>
      <w:sz w:val="19"/>
      <w:szCs w:val="19"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:highlight w:val="cyan"/>
    </w:rPr>
    <w:t>fff</w:t>
  </w:r>
  <w:bookmarkStart w:id="1" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="1"/>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="2"/>
    </w:numPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>ff</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>

</pkg:xmlData>
'@

$insertAt.InsertXML($newXml)

Write-Output "done"
